# Auto-generated BRVM data refresh (GitHub Actions style update)
$wb = $excel.ActiveWorkbook

# --- Sheet "Recommandations": refresh rows 2-43 (A:G) with latest data ---
$wsReco = $wb.Worksheets.Item("Recommandations")

$sheet1Data = @(
    @("BRVM - SERVICES PUBLICS", 0, 8, 3339.68, 105.84, "🟡 Observer", "➖ Neutre"),
    @("AIR LIQUIDE CI", 0, 4, 2740, 680, "🟡 Observer", "➖ Neutre"),
    @("NEI-CEDA CI", 0, 4, 2680, 685, "🟡 Observer", "➖ Neutre"),
    @("BRVM - AUTRES SECTEURS", 0, 4, 2313.36, 564.23, "🟡 Observer", "➖ Neutre"),
    @("BRVM - DISTRIBUTION", 0, 4, 1722.95, 444.94, "🟡 Observer", "➖ Neutre"),
    @("BRVM - TRANSPORT", 0, 4, 1445.38, 359.82, "🟡 Observer", "➖ Neutre"),
    @("BRVM - AGRICULTURE", 0, 4, 1443.82, 364.5, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DISCRETIONNAIRE", 0, 4, 551.76, 143.68, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRESTIGE", 0, 4, 547.38, 137.37, "🟡 Observer", "➖ Neutre"),
    @("BRVM - FINANCES", 0, 4, 542.6900000000001, 135.6, "🟡 Observer", "➖ Neutre"),
    @("BRVM - SERVICES FINANCIERS", 0, 4, 533.35, 133.27, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIELS", 0, 4, 487.72, 124.43, "🟡 Observer", "➖ Neutre"),
    @("BRVM - ENERGIE", 0, 4, 431.23, 108.58, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIE                  (**)", 0, 2, 414.76, 207.68, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRINCIPAL                    (**)", 0, 2, 387.23, 193.83, "🟡 Observer", "➖ Neutre"),
    @("BRVM - TELECOMMUNICATIONS", 0, 4, 381.09, 95.94, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DE BASE         (**)", 0, 2, 372.63, 187.07, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIE                (**)", 0, 1, 206.88, 206.88, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRINCIPAL                  (**)", 0, 1, 193.58, 193.58, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DE BASE              (**)", 0, 1, 187.55, 187.55, "🟡 Observer", "➖ Neutre"),
    @("FILTISAC CI (FTSC)", 4, 0, 29.51, 7.39, "🟢 Achat", "✅ Renforcer"),
    @("ORAGROUP TOGO (ORGT)", 2, 0, 14.72, 7.41, "🟡 Observer", "➖ Neutre"),
    @("BERNABE CI (BNBC)", 2, 0, 12.96, 6.51, "🟡 Observer", "➖ Neutre"),
    @("CFAO MOTORS CI (CFAC)", 2, 1, 8.85, 7.46, "🟡 Observer", "👀 À surveiller"),
    @("TRACTAFRIC MOTORS CI (PRSC)", 2, 1, 8.19, 7.34, "🟡 Observer", "👀 À surveiller"),
    @("VIVO ENERGY CI (SHEC)", 2, 0, 4.19, 2.92, "🟡 Observer", "➖ Neutre"),
    @("ONATEL BF (ONTBF)", 1, 0, 2.61, 2.61, "🟡 Observer", "➖ Neutre"),
    @("BICI CI (BICC)", 1, 0, 2.2, 2.2, "🟡 Observer", "➖ Neutre"),
    @("NSIA BANQUE COTE D'IVOIRE (NSBC)", 1, 1, 0.49, -2.01, "🟡 Observer", "👀 À surveiller"),
    @("SICABLE CI (CABC)", 1, 1, 0.25, 5.2, "🟡 Observer", "👀 À surveiller"),
    @("TOTAL", 0, 4, 0, 0, "🟡 Observer", "➖ Neutre"),
    @("UNIWAX CI (UNXC)", 1, 1, -0.25, 7.23, "🟡 Observer", "👀 À surveiller"),
    @("SAPH CI (SPHC)", 1, 1, -0.32, 5.94, "🟡 Observer", "👀 À surveiller"),
    @("AFRICA GLOBAL LOGISTICS CI (SDSC)", 0, 1, -1.68, -1.68, "🟡 Observer", "➖ Neutre"),
    @("TOTALENERGIES MARKETING CI (TTLC)", 0, 1, -2.04, -2.04, "🟡 Observer", "➖ Neutre"),
    @("SOLIBRA CI (SLBC)", 0, 1, -2.75, -2.75, "🟡 Observer", "➖ Neutre"),
    @("NEI-CEDA CI (NEIC)", 0, 1, -2.94, -2.94, "🟡 Observer", "➖ Neutre"),
    @("SETAO CI (STAC)", 0, 2, -4.81, -2.91, "🟡 Observer", "➖ Neutre"),
    @("AIR LIQUIDE CI (SIVC)", 0, 1, -7.35, -7.35, "🟡 Observer", "➖ Neutre"),
    @("UNILEVER CI (UNLC)", 0, 1, -7.5, -7.5, "🟡 Observer", "➖ Neutre"),
    @("SUCRIVOIRE (SCRC)", 0, 3, -13.69, -3.43, "🔴 Vente", "⚠️ Risque de décrochage"),
    @("LOTERIE NATIONALE DU BENIN (LNBB)", 0, 2, -14.92, -7.45, "🟡 Observer", "➖ Neutre")
)

for ($i = 0; $i -lt $sheet1Data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $sheet1Data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $wsReco.Cells.Item($rowNum, $c + 1).Value = $rowVals[$c]
    }
}

# Row 44 (LOTERIE NATIONALE DU BENIN) no longer appears in the refreshed data set
$wsReco.Rows.Item(44).Delete()

# --- Sheet "Top_YTD": refresh column B (rows 2-11) with latest progression values ---
$wsYtd = $wb.Worksheets.Item("Top_YTD")

$sheet2Data = @(
    8342409.89
    379602.44
    351237.74
    211315.75
    79203.55
    45198.18
    45039.76
    3103.91
    3046.68
    2984.85
)

for ($i = 0; $i -lt $sheet2Data.Count; $i++) {
    $rowNum = $i + 2
    $wsYtd.Cells.Item($rowNum, 2).Value = $sheet2Data[$i]
}
